$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.443.50"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "2.200.82"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.33"
$ws.Range("E5").Value = "  -2.09%  "

$ws.Range("E6").Value = "  -1.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.41"
$ws.Range("E7").Value = "  -3.43%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("E9").Value = "  -4.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.21"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  -4.08%  "

$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.82"
$ws.Range("E13").Value = "  -4.29%  "

$ws.Range("D14").Value = "2.530.42"
$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.08"
$ws.Range("E15").Value = "  -2.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.821"
$ws.Range("E16").Value = "  -3.67%  "

$ws.Range("D17").Value = "2.198.34"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "41.378.93"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("E19").Value = "  -10.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.08"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.17"
$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("E22").Value = "  +6.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "227.50"
$ws.Range("E23").Value = "  -1.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -5.97%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  -5.77%  "

$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  -2.90%  "

$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.28"
$ws.Range("E30").Value = "  -0.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.28"
$ws.Range("E31").Value = "  -3.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0786"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.39"
$ws.Range("E33").Value = "  +3.89%  "

$ws.Range("E34").Value = "  -7.01%  "

$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("E36").Value = "  -9.50%  "

$ws.Range("E37").Value = "  -4.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0297"
$ws.Range("E38").Value = "  -2.17%  "

$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.08"
$ws.Range("E40").Value = "  -3.85%  "

$ws.Range("E41").Value = "  -1.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.92"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.193"
$ws.Range("E43").Value = "  -3.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.54"
$ws.Range("E44").Value = "  -2.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("E45").Value = "  -3.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.65"
$ws.Range("E46").Value = "  -3.77%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("E48").Value = "  -3.25%  "

$ws.Range("E49").Value = "  -2.43%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.67"
$ws.Range("E50").Value = "  -1.63%  "

$ws.Range("D51").Value = "2.409.23"
$ws.Range("E51").Value = "  -1.14%  "
